# Re-apply the Design > Themes switch (Integral -> Office Theme) together
# with the table style refresh that PowerPoint performs on the three
# tables that were still carrying the old theme's table style id.

$p = $ppt.ActivePresentation

# --- 1. Tables on slides 14, 15 and 16 pick up the new built-in table
#        style id that ships with the "Office Theme" design. ---
$newTableStyleId = "{922E8737-DA86-49AD-8E98-EBA72BDA6DB2}"
foreach ($slideIndex in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIndex)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shp = $slide.Shapes.Item($i)
        if ($shp.HasTable) {
            $shp.Table.ApplyStyle($newTableStyleId)
        }
    }
}

# --- 2. Switch the presentation's design theme from "Integral" (Red
#        Violet color scheme) back to the default "Office Theme" color
#        scheme. The color scheme is expressed as 12 theme color slots
#        in the fixed order: dk1, lt1, dk2, lt2, accent1-6, hlink,
#        folHlink. ---
$officeThemeColors = @(
    0,            # dk1      000000
    16777215,     # lt1      FFFFFF
    6968388,      # dk2      44546A
    15132391,     # lt2      E7E6E6
    13998939,     # accent1  5B9BD5
    3243501,      # accent2  ED7D31
    10855845,     # accent3  A5A5A5
    49407,        # accent4  FFC000
    12874308,     # accent5  4472C4
    4697456,      # accent6  70AD47
    12673797,     # hlink    0563C1
    7491477       # folHlink 954F72
)

$scheme = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le $scheme.Count; $i++) {
    $scheme.Item($i).RGB = $officeThemeColors[$i - 1]
}
